# Apply updated FFXIV market price data to the Halicarnassus Profits workbook.
# Generated from the authoritative commit diff: for each sheet/row, refresh the
# price/profit columns (H-N) with the scraper run's latest values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1245.5
$ws.Range("I6").Value = 271
$ws.Range("K6").Value = 813
$ws.Range("M6").Value = -701

$ws.Range("H28").Value = 6362.625
$ws.Range("I28").Value = 829.7
$ws.Range("J28").Value = 15584.167
$ws.Range("K28").Value = 829.7
$ws.Range("L28").Value = 15584.167
$ws.Range("M28").Value = -344.7
$ws.Range("N28").Value = -16554.167

$ws.Range("H31").Value = 53
$ws.Range("I31").Value = 53
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 159
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 71
$ws.Range("N31").ClearContents()

$ws.Range("H33").Value = 354
$ws.Range("I33").Value = 354
$ws.Range("K33").Value = 354
$ws.Range("M33").Value = -125

$ws.Range("H38").Value = 387.0909
$ws.Range("I38").Value = 263.33334
$ws.Range("J38").Value = 944
$ws.Range("K38").Value = 790.0000200000001
$ws.Range("L38").Value = 2832
$ws.Range("M38").Value = -418.0000200000001
$ws.Range("N38").Value = -3576

$ws.Range("H40").Value = 5280.7896
$ws.Range("I40").Value = 4427.1816
$ws.Range("J40").Value = 6454.5
$ws.Range("K40").Value = 4427.1816
$ws.Range("L40").Value = 6454.5
$ws.Range("M40").Value = -4252.1816
$ws.Range("N40").Value = -6804.5

$ws.Range("H62").Value = 7701.769
$ws.Range("I62").Value = 4015.375
$ws.Range("J62").Value = 13600
$ws.Range("K62").Value = 4015.375
$ws.Range("L62").Value = 13600
$ws.Range("M62").Value = -3391.375
$ws.Range("N62").Value = -14848

$ws.Range("H65").Value = 7701.769
$ws.Range("I65").Value = 4015.375
$ws.Range("J65").Value = 13600
$ws.Range("K65").Value = 20076.875
$ws.Range("L65").Value = 68000
$ws.Range("M65").Value = -16956.875
$ws.Range("N65").Value = -74240

$ws.Range("H82").Value = 350
$ws.Range("I82").Value = 350
$ws.Range("K82").Value = 1050
$ws.Range("M82").Value = -644

$ws.Range("H85").Value = 350
$ws.Range("I85").Value = 350
$ws.Range("K85").Value = 1050
$ws.Range("M85").Value = 354

$ws.Range("H116").Value = 6612.375
$ws.Range("I116").Value = 8000
$ws.Range("K116").Value = 8000
$ws.Range("M116").Value = -4558

$ws.Range("H138").Value = 2645.3
$ws.Range("J138").Value = 4733.3335
$ws.Range("L138").Value = 14200.0005
$ws.Range("N138").Value = -24480.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H45").Value = 2989.125
$ws.Range("I45").Value = 2152.1667
$ws.Range("J45").Value = 5500
$ws.Range("K45").Value = 2152.1667
$ws.Range("L45").Value = 5500
$ws.Range("M45").Value = -1775.1667
$ws.Range("N45").Value = -6254

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1085.5883
$ws.Range("I105").Value = 893.2143
$ws.Range("K105").Value = 893.2143
$ws.Range("M105").Value = 853.7857

$ws.Range("H132").Value = 130780
$ws.Range("J132").Value = 130780
$ws.Range("L132").Value = 130780
$ws.Range("N132").Value = -140900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 16426.25
$ws.Range("I6").Value = 2901
$ws.Range("J6").Value = 57002
$ws.Range("K6").Value = 2901
$ws.Range("L6").Value = 57002
$ws.Range("M6").Value = -2788
$ws.Range("N6").Value = -57228

$ws.Range("H19").Value = 9571.166999999999
$ws.Range("I19").Value = 84.2
$ws.Range("K19").Value = 84.2
$ws.Range("M19").Value = 85.8

$ws.Range("H22").Value = 1140.3572
$ws.Range("I22").Value = 773.1111
$ws.Range("K22").Value = 773.1111
$ws.Range("M22").Value = -423.1111

$ws.Range("H24").Value = 9571.166999999999
$ws.Range("I24").Value = 84.2
$ws.Range("K24").Value = 84.2
$ws.Range("M24").Value = 85.8

$ws.Range("H106").Value = 24091.223
$ws.Range("J106").Value = 24091.223
$ws.Range("L106").Value = 24091.223
$ws.Range("N106").Value = -26615.223

$ws.Range("H122").Value = 871.7143
$ws.Range("J122").Value = 1221
$ws.Range("L122").Value = 3663
$ws.Range("N122").Value = -8563

$ws.Range("H125").Value = 63270.668
$ws.Range("J125").Value = 63270.668
$ws.Range("L125").Value = 63270.668
$ws.Range("N125").Value = -68190.66800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 98
$ws.Range("I6").Value = 98
$ws.Range("K6").Value = 294
$ws.Range("M6").Value = -181

$ws.Range("H17").Value = 392.95456
$ws.Range("I17").Value = 77.07692
$ws.Range("K17").Value = 231.23076
$ws.Range("M17").Value = -62.23076

$ws.Range("H23").Value = 156.33333
$ws.Range("I23").Value = 135.5
$ws.Range("J23").Value = 166.75
$ws.Range("K23").Value = 406.5
$ws.Range("L23").Value = 500.25
$ws.Range("M23").Value = -171.5
$ws.Range("N23").Value = -970.25

$ws.Range("H132").Value = 1309.6666
$ws.Range("I132").Value = 972.25
$ws.Range("K132").Value = 8750.25
$ws.Range("M132").Value = -6220.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 38.714287
$ws.Range("I2").Value = 37.625
$ws.Range("J2").Value = 40.166668
$ws.Range("K2").Value = 37.625
$ws.Range("L2").Value = 40.166668
$ws.Range("M2").Value = 75.375
$ws.Range("N2").Value = -266.166668

$ws.Range("H21").Value = 20000
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20346

$ws.Range("H30").Value = 20000
$ws.Range("J30").Value = 20000
$ws.Range("L30").Value = 20000
$ws.Range("N30").Value = -20210

$ws.Range("H92").Value = 5343.75
$ws.Range("J92").Value = 5343.75
$ws.Range("L92").Value = 5343.75
$ws.Range("N92").Value = -9087.75

$ws.Range("H97").Value = 1224.4445
$ws.Range("I97").Value = 1168.1666
$ws.Range("J97").Value = 1337
$ws.Range("K97").Value = 1168.1666
$ws.Range("L97").Value = 1337
$ws.Range("M97").Value = -672.1666
$ws.Range("N97").Value = -2329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4975.533
$ws.Range("I7").Value = 3467.75
$ws.Range("J7").Value = 6698.7144
$ws.Range("K7").Value = 3467.75
$ws.Range("L7").Value = 6698.7144
$ws.Range("M7").Value = -3355.75
$ws.Range("N7").Value = -6922.7144

$ws.Range("H16").Value = 307.1
$ws.Range("I16").Value = 307.1
$ws.Range("K16").Value = 307.1
$ws.Range("M16").Value = -137.1

$ws.Range("H22").Value = 866.3333
$ws.Range("I22").Value = 876
$ws.Range("K22").Value = 876
$ws.Range("M22").Value = -581

$ws.Range("H27").Value = 866.3333
$ws.Range("I27").Value = 876
$ws.Range("K27").Value = 876
$ws.Range("M27").Value = -769

$ws.Range("H94").Value = 145000
$ws.Range("J94").Value = 145000
$ws.Range("L94").Value = 145000
$ws.Range("N94").Value = -146352

$ws.Range("H126").Value = 4975.533
$ws.Range("I126").Value = 3467.75
$ws.Range("J126").Value = 6698.7144
$ws.Range("K126").Value = 10403.25
$ws.Range("L126").Value = 20096.1432
$ws.Range("M126").Value = -7933.25
$ws.Range("N126").Value = -25036.1432

$ws.Range("H132").Value = 3302.1333
$ws.Range("I132").Value = 3000.2727
$ws.Range("K132").Value = 9000.8181
$ws.Range("M132").Value = -6470.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8916.666999999999
$ws.Range("J62").Value = 10875
$ws.Range("L62").Value = 10875
$ws.Range("N62").Value = -12123

$ws.Range("H65").Value = 8916.666999999999
$ws.Range("J65").Value = 10875
$ws.Range("L65").Value = 54375
$ws.Range("N65").Value = -60615

$ws.Range("H107").Value = 849.5
$ws.Range("I107").Value = 849.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2548.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -628.5
$ws.Range("N107").ClearContents()

$ws.Range("H132").Value = 3939.6428
$ws.Range("I132").Value = 2785
$ws.Range("K132").Value = 8355
$ws.Range("M132").Value = -5825
